$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "528.0 - 572.0"
$ws.Range("B2").Value = 564.5
$ws.Range("C2").Value = 555.9
$ws.Range("D2").Value = 546
$ws.Range("E2").Value = 551.2
$ws.Range("F2").Value = 560.2

$ws.Range("A3").Value = "336.0 - 364.0"
$ws.Range("B3").Value = 356.7
$ws.Range("C3").Value = 357
$ws.Range("D3").Value = 343.6
$ws.Range("E3").Value = 342.4
$ws.Range("F3").Value = 355.2

$ws.Range("A4").Value = "192.0 - 208.0"
$ws.Range("B4").Value = 197.4
$ws.Range("C4").Value = 200.4
$ws.Range("D4").Value = 199.5
$ws.Range("E4").Value = 205.3
$ws.Range("F4").Value = 205.5
